$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet and relocate the selection ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "2 Switch"
$ws1.Range("G22").Select()

# --- Add the new "3 Switch" sheet right after "2 Switch" ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "3 Switch"

# --- Headers ---
$ws2.Range("A1").Value = "R1"
$ws2.Range("B1").Value = "R2"
$ws2.Range("C1").Value = "R3"
$ws2.Range("D1").Value = "S1 VDC"
$ws2.Range("E1").Value = "S1 A2D"
$ws2.Range("F1").Value = "S2 VDC"
$ws2.Range("G1").Value = "S2 A2D"
$ws2.Range("H1").Value = "S3 VDC"
$ws2.Range("I1").Value = "S3 A2D"

# --- Values / formulas ---
$ws2.Range("A2").Value = 2200
$ws2.Range("B2").Value = 1000
$ws2.Range("C2").Value = 620
$ws2.Range("D2").Value = 0
$ws2.Range("E2").Value = 0
$ws2.Range("F2").Formula = "=(B2/(A2+B2))*5"
$ws2.Range("G2").Formula = "=(B2/(A2+B2))*1024"
$ws2.Range("H2").Formula = "=(C2/(A2+B2+C2))*5"
$ws2.Range("I2").Formula = "=(C2/(A2+B2+C2))*1024"

# --- Number formats (reuse the workbook's existing "0.00"/"0" styles) ---
$ws2.Range("F1:F2").NumberFormat = "0.00"
$ws2.Range("H1:H2").NumberFormat = "0.00"
$ws2.Range("I1:I2").NumberFormat = "0"

# --- Match the original's portrait page setup ---
$ws2.PageSetup.Orientation = 1

# --- Selection for the new, now-active sheet ---
$ws2.Range("C3").Select()
